$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = $null
$ws.Range("A2").Value = 44
$ws.Range("A3").Value = 67

$ws.Range("A3").Select()
